$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# --- 1) "Total number of PhD awarded in last five years: ..." paragraph ---
# Split "Total number of PhD awarded in last " into "Total number of " + "PG"
# (bold run) + the _GoBack bookmark + " awarded in last ", keeping "five"
# and " years: _____________" runs untouched.
$rng1 = $d.Content
$rng1.Find.Execute("Total number of PhD awarded in last five years: _____________") | Out-Null
$para1 = $rng1.Paragraphs(1).Range

$xml1 = '<w:p xmlns:w="' + $wNs + '" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="6AD66ECF" w14:textId="77777777" w:rsidR="00EF76DE" w:rsidRPr="00BD0FFD" w:rsidRDefault="00EF76DE" w:rsidP="00EF76DE">' +
  '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
  '<w:r w:rsidRPr="00BD0FFD"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Total number of </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t>PG</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> awarded in last </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t>five</w:t></w:r>' +
  '<w:r w:rsidRPr="00BD0FFD"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> years: _____________</w:t></w:r>' +
  '</w:p>'

$para1.InsertXML($xml1)

# --- 2) "Annexure-VII" heading paragraph ---
# Remove the (now duplicate) _GoBack bookmark that used to sit here, and add
# a lastRenderedPageBreak marker inside the "Annexure-" run.
$rng2 = $d.Content
$rng2.Find.Execute("Annexure-VII") | Out-Null
$para2 = $rng2.Paragraphs(1).Range

$xml2 = '<w:p xmlns:w="' + $wNs + '" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="01450B5F" w14:textId="23C99885" w:rsidR="00C2072F" w:rsidRDefault="00C2072F" w:rsidP="00C2072F">' +
  '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="right"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:color w:val="333399"/><w:sz w:val="28"/><w:szCs w:val="24"/><w:u w:val="single"/><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
  '<w:r w:rsidRPr="00813100"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:color w:val="333399"/><w:sz w:val="28"/><w:szCs w:val="24"/><w:u w:val="single"/><w:lang w:val="en-GB"/></w:rPr><w:lastRenderedPageBreak/><w:t>Annexure-</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:color w:val="333399"/><w:sz w:val="28"/><w:szCs w:val="24"/><w:u w:val="single"/><w:lang w:val="en-GB"/></w:rPr><w:t>V</w:t></w:r>' +
  '<w:r w:rsidRPr="00813100"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:color w:val="333399"/><w:sz w:val="28"/><w:szCs w:val="24"/><w:u w:val="single"/><w:lang w:val="en-GB"/></w:rPr><w:t>I</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:color w:val="333399"/><w:sz w:val="28"/><w:szCs w:val="24"/><w:u w:val="single"/><w:lang w:val="en-GB"/></w:rPr><w:t>I</w:t></w:r>' +
  '</w:p>'

$para2.InsertXML($xml2)
